# Apply updated res_bus/vm_pu.xlsx voltage-magnitude results for the 380 kV case.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.038825164163572
$rowBF[0,2] = 1.039652455884032
$rowBF[0,3] = 1.052365701346886
$rowBF[0,4] = 1.059753805040696
$ws.Range("B2:F2").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.037340896980113
$rowIN[0,1] = 1.043920502206693
$rowIN[0,2] = 1.042437034650083
$rowIN[0,3] = 1.055114634199702
$rowIN[0,4] = 1.062482448437402
$rowIN[0,5] = 1.045402989059181
$ws.Range("I2:N2").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.039910863128316
$rowBF[0,2] = 1.040440374574999
$rowBF[0,3] = 1.053437375232783
$rowBF[0,4] = 1.060895687153399
$ws.Range("B3:F3").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.037574510516013
$rowIN[0,1] = 1.044650471684014
$rowIN[0,2] = 1.043035274737793
$rowIN[0,3] = 1.05599847495508
$rowIN[0,4] = 1.06343779371873
$rowIN[0,5] = 1.04613399517689
$ws.Range("I3:N3").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.04061338324055
$rowBF[0,2] = 1.040949850088233
$rowBF[0,3] = 1.054131151800219
$rowBF[0,4] = 1.061634891125991
$ws.Range("B4:F4").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.037723872525405
$rowIN[0,1] = 1.04512225153069
$rowIN[0,2] = 1.043421361235424
$rowIN[0,3] = 1.056570113134937
$rowIN[0,4] = 1.064055710786332
$rowIN[0,5] = 1.046606445005062
$ws.Range("I4:N4").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.040908723237402
$rowBF[0,2] = 1.041163946807452
$rowBF[0,3] = 1.054422894971772
$rowBF[0,4] = 1.061945731605456
$ws.Range("B5:F5").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.037786232780766
$rowIN[0,1] = 1.045320453885511
$rowIN[0,2] = 1.043583428808287
$rowIN[0,3] = 1.056810366197605
$rowIN[0,4] = 1.064315422089358
$rowIN[0,5] = 1.04680492882996
$ws.Range("I5:N5").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.040958312167981
$rowBF[0,2] = 1.041199889485813
$rowBF[0,3] = 1.054471884597342
$rowBF[0,4] = 1.061997927734333
$ws.Range("B6:F6").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.037796678029844
$rowIN[0,1] = 1.045353725080999
$rowIN[0,2] = 1.04361062639915
$rowIN[0,3] = 1.056850702013534
$rowIN[0,4] = 1.064359025182842
$rowIN[0,5] = 1.046838247274363
$ws.Range("I6:N6").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.040617329585318
$rowBF[0,2] = 1.040952711200684
$rowBF[0,3] = 1.054135049774312
$rowBF[0,4] = 1.061639044280579
$ws.Range("B7:F7").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.037724707482122
$rowIN[0,1] = 1.045124900446389
$rowIN[0,2] = 1.043423527746145
$rowIN[0,3] = 1.056573323658107
$rowIN[0,4] = 1.064059181301213
$rowIN[0,5] = 1.046609097682524
$ws.Range("I7:N7").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.039192081438935
$rowBF[0,2] = 1.039918811074172
$rowBF[0,3] = 1.052727809859556
$rowBF[0,4] = 1.060139641234231
$ws.Range("B8:F8").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.03742022065588
$rowIN[0,1] = 1.044167314591
$rowIN[0,2] = 1.042639422678679
$rowIN[0,3] = 1.055413387392443
$rowIN[0,4] = 1.062805364732367
$rowIN[0,5] = 1.045650151945384
$ws.Range("I8:N8").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.036680592517538
$rowBF[0,2] = 1.038094205336963
$rowBF[0,3] = 1.050250612578387
$rowBF[0,4] = 1.057500020470951
$ws.Range("B9:F9").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.03686988799859
$rowIN[0,1] = 1.042475643055101
$rowIN[0,2] = 1.041249966647628
$rowIN[0,3] = 1.053367393215216
$rowIN[0,4] = 1.06059402144165
$rowIN[0,5] = 1.043956078041866
$ws.Range("I9:N9").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.035006219901886
$rowBF[0,2] = 1.036875982406863
$rowBF[0,3] = 1.048600847700992
$rowBF[0,4] = 1.055741952748326
$ws.Range("B10:F10").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.036493741147862
$rowIN[0,1] = 1.041344979654882
$rowIN[0,2] = 1.040318454116175
$rowIN[0,3] = 1.052002020910893
$rowIN[0,4] = 1.059118467450013
$rowIN[0,5] = 1.042823808969931
$ws.Range("I10:N10").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.034281178587605
$rowBF[0,2] = 1.03634805122602
$rowBF[0,3] = 1.047886881815552
$rowBF[0,4] = 1.054981084099288
$ws.Range("B11:F11").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.036328671429852
$rowIN[0,1] = 1.040854704755316
$rowIN[0,2] = 1.039913864976925
$rowIN[0,3] = 1.051410470953335
$rowIN[0,4] = 1.058479219051003
$rowIN[0,5] = 1.042332837823772
$ws.Range("I11:N11").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.03401186107775
$rowBF[0,2] = 1.03615188932918
$rowBF[0,3] = 1.047621741819927
$rowBF[0,4] = 1.054698521149793
$ws.Range("B12:F12").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.036267027400159
$rowIN[0,1] = 1.040672490902056
$rowIN[0,2] = 1.039763396567115
$rowIN[0,3] = 1.051190692362053
$rowIN[0,4] = 1.058241724716691
$rowIN[0,5] = 1.04215036520594
$ws.Range("I12:N12").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.034069630832295
$rowBF[0,2] = 1.036193969648155
$rowBF[0,3] = 1.04767861259731
$rowBF[0,4] = 1.054759129265528
$ws.Range("B13:F13").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.036280265176001
$rowIN[0,1] = 1.040711581098779
$rowIN[0,2] = 1.039795680967989
$rowIN[0,3] = 1.051237837910592
$rowIN[0,4] = 1.058292670281694
$rowIN[0,5] = 1.042189510915226
$ws.Range("I13:N13").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.034258916822697
$rowBF[0,2] = 1.036331837741924
$rowBF[0,3] = 1.047864964064891
$rowBF[0,4] = 1.054957726186688
$ws.Range("B14:F14").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.036323582642197
$rowIN[0,1] = 1.040839645022785
$rowIN[0,2] = 1.039901431005277
$rowIN[0,3] = 1.051392305020398
$rowIN[0,4] = 1.058459588697322
$rowIN[0,5] = 1.042317756704694
$ws.Range("I14:N14").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.034375541496429
$rowBF[0,2] = 1.036416774251913
$rowBF[0,3] = 1.047979789125383
$rowBF[0,4] = 1.055080095914851
$ws.Range("B15:F15").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.036350228265321
$rowIN[0,1] = 1.04091853565476
$rowIN[0,2] = 1.039966562448477
$rowIN[0,3] = 1.051487470590795
$rowIN[0,4] = 1.058562426135876
$rowIN[0,5] = 1.042396759370414
$ws.Range("I15:N15").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.03505433786729
$rowBF[0,2] = 1.036911010367046
$rowBF[0,3] = 1.048648239519812
$rowBF[0,4] = 1.055792457191417
$ws.Range("B16:F16").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.036504650028908
$rowIN[0,1] = 1.041377503024011
$rowIN[0,2] = 1.040345279281738
$rowIN[0,3] = 1.052041273054608
$rowIN[0,4] = 1.059160885407559
$rowIN[0,5] = 1.042856378525974
$ws.Range("I16:N16").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.035480121098227
$rowBF[0,2] = 1.037220915977043
$rowBF[0,3] = 1.04906764581584
$rowBF[0,4] = 1.056239405673055
$ws.Range("B17:F17").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.036600926944817
$rowIN[0,1] = 1.041665216183119
$rowIN[0,2] = 1.040582506893504
$rowIN[0,3] = 1.052388568873902
$rowIN[0,4] = 1.059536196150685
$rowIN[0,5] = 1.043144500270769
$ws.Range("I17:N17").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.035728470711877
$rowBF[0,2] = 1.037401636919231
$rowBF[0,3] = 1.049312316373616
$rowBF[0,4] = 1.056500140463716
$ws.Range("B18:F18").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.036656871787165
$rowIN[0,1] = 1.041832967824022
$rowIN[0,2] = 1.040720758289647
$rowIN[0,3] = 1.052591108372816
$rowIN[0,4] = 1.059755077342124
$rowIN[0,5] = 1.043312490138243
$ws.Range("I18:N18").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.035813151106631
$rowBF[0,2] = 1.037463251000007
$rowBF[0,3] = 1.049395749160389
$rowBF[0,4] = 1.056589050631416
$ws.Range("B19:F19").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.036675911595945
$rowIN[0,1] = 1.041890155512895
$rowIN[0,2] = 1.040767878180575
$rowIN[0,3] = 1.052660163642463
$rowIN[0,4] = 1.059829704888601
$rowIN[0,5] = 1.043369759040194
$ws.Range("I19:N19").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.035434438855843
$rowBF[0,2] = 1.03718767035523
$rowBF[0,3] = 1.049022643565362
$rowBF[0,4] = 1.056191448475772
$ws.Range("B20:F20").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.036590619255246
$rowIN[0,1] = 1.041634354171971
$rowIN[0,2] = 1.040557066971545
$rowIN[0,3] = 1.052351310659867
$rowIN[0,4] = 1.059495932091115
$rowIN[0,5] = 1.043113594432026
$ws.Range("I20:N20").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.034203176957897
$rowBF[0,2] = 1.036291240810066
$rowBF[0,3] = 1.047810086584142
$rowBF[0,4] = 1.054899242758742
$ws.Range("B21:F21").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.036310835831125
$rowIN[0,1] = 1.040801936259834
$rowIN[0,2] = 1.039870295389151
$rowIN[0,3] = 1.051346819705798
$rowIN[0,4] = 1.058410436740346
$rowIN[0,5] = 1.042279994390974
$ws.Range("I21:N21").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.033429005064812
$rowBF[0,2] = 1.035727245468119
$rowBF[0,3] = 1.047048043464904
$rowBF[0,4] = 1.054087114426811
$ws.Range("B22:F22").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.036133016883278
$rowIN[0,1] = 1.040277960677806
$rowIN[0,2] = 1.039437418766028
$rowIN[0,3] = 1.050714963474243
$rowIN[0,4] = 1.057727658950062
$rowIN[0,5] = 1.041755274703519
$ws.Range("I22:N22").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.033839411231158
$rowBF[0,2] = 1.036026265488982
$rowBF[0,3] = 1.04745198486999
$rowBF[0,4] = 1.054517607650708
$ws.Range("B23:F23").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.036227462901364
$rowIN[0,1] = 1.040555787132644
$rowIN[0,2] = 1.039666996849299
$rowIN[0,3] = 1.051049950284913
$rowIN[0,4] = 1.058089639399424
$rowIN[0,5] = 1.04203349570379
$ws.Range("I23:N23").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.035455080702347
$rowBF[0,2] = 1.037202692749612
$rowBF[0,3] = 1.049042978027741
$rowBF[0,4] = 1.056213118153414
$ws.Range("B24:F24").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.036595277510997
$rowIN[0,1] = 1.041648299591909
$rowIN[0,2] = 1.040568562546206
$rowIN[0,3] = 1.052368146142557
$rowIN[0,4] = 1.059514125782714
$rowIN[0,5] = 1.043127559656059
$ws.Range("I24:N24").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.037329878366224
$rowBF[0,2] = 1.038566232029918
$rowBF[0,3] = 1.050890727272617
$rowBF[0,4] = 1.058182128641149
$ws.Range("B25:F25").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.037013794891748
$rowIN[0,1] = 1.042913488514703
$rowIN[0,2] = 1.041610093134266
$rowIN[0,3] = 1.053896573876167
$rowIN[0,4] = 1.0611708694453141
$rowIN[0,5] = 1.044394545292243
$ws.Range("I25:N25").Value = $rowIN

Write-Host "Updated vm_pu results for rows 2-25 (24 buses x 11 columns)."